$wb = $excel.ActiveWorkbook

# Sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2798.5
$ws.Range("I40").Value = 2888.9
$ws.Range("J40").Value = 2572.5
$ws.Range("K40").Value = 2888.9
$ws.Range("L40").Value = 2572.5
$ws.Range("M40").Value = -2713.9
$ws.Range("N40").Value = -2922.5

# Sheet ALC, row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 826.2381
$ws.Range("I98").Value = 855.05
$ws.Range("J98").Value = 250
$ws.Range("K98").Value = 855.05
$ws.Range("L98").Value = 250
$ws.Range("M98").Value = 642.95
$ws.Range("N98").Value = -3246

# Sheet ALC, row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1071.1305
$ws.Range("J112").Value = 1219.7778
$ws.Range("L112").Value = 3659.3334
$ws.Range("N112").Value = -5875.3334

# Sheet ALC, row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 826.2381
$ws.Range("I122").Value = 855.05
$ws.Range("J122").Value = 250
$ws.Range("K122").Value = 2565.15
$ws.Range("L122").Value = 750
$ws.Range("M122").Value = -115.1499999999996
$ws.Range("N122").Value = -5650

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5372.4688
$ws.Range("I32").Value = 3661.93
$ws.Range("K32").Value = 3661.93
$ws.Range("M32").Value = -3374.93

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1107.65
$ws.Range("I45").Value = 814.5
$ws.Range("K45").Value = 814.5
$ws.Range("M45").Value = -437.5

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 931.2222
$ws.Range("I61").Value = 810.82355
$ws.Range("J61").Value = 1135.9
$ws.Range("K61").Value = 810.82355
$ws.Range("L61").Value = 1135.9
$ws.Range("M61").Value = -598.82355
$ws.Range("N61").Value = -1559.9

# Sheet ARM, row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1750
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 1750
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 1750
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -2562

# Sheet ARM, row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1750
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1750
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 1750
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -4558

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2066
$ws.Range("I122").Value = 1826.1052
$ws.Range("J122").Value = 2480.3635
$ws.Range("K122").Value = 5478.3156
$ws.Range("L122").Value = 7441.0905
$ws.Range("M122").Value = -3028.3156
$ws.Range("N122").Value = -12341.0905

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 931.2222
$ws.Range("I136").Value = 810.82355
$ws.Range("J136").Value = 1135.9
$ws.Range("K136").Value = 2432.47065
$ws.Range("L136").Value = 3407.7
$ws.Range("M136").Value = 117.5293500000002
$ws.Range("N136").Value = -8507.700000000001

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1342.7587
$ws.Range("I99").Value = 791.0526
$ws.Range("J99").Value = 2391
$ws.Range("K99").Value = 791.0526
$ws.Range("L99").Value = 2391
$ws.Range("M99").Value = 706.9474
$ws.Range("N99").Value = -5387

# Sheet CRP, row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1560.1
$ws.Range("I16").Value = 1509.8889
$ws.Range("K16").Value = 1509.8889
$ws.Range("M16").Value = -1222.8889

# Sheet CRP, row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1390
$ws.Range("I105").Value = 1390
$ws.Range("K105").Value = 1390
$ws.Range("M105").Value = 357

# Sheet CRP, row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1560.1
$ws.Range("I113").Value = 1509.8889
$ws.Range("K113").Value = 1509.8889
$ws.Range("M113").Value = 660.1111000000001

# Sheet CRP, row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 423926.88
$ws.Range("I122").Value = 59709.59
$ws.Range("J122").Value = 1111892.9
$ws.Range("K122").Value = 179128.77
$ws.Range("L122").Value = 3335678.7
$ws.Range("M122").Value = -176678.77
$ws.Range("N122").Value = -3340578.7

# Sheet CUL, row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 91618.63
$ws.Range("I92").Value = 143536.58
$ws.Range("J92").Value = 762.25
$ws.Range("K92").Value = 430609.74
$ws.Range("L92").Value = 2286.75
$ws.Range("M92").Value = -429361.74
$ws.Range("N92").Value = -4782.75

# Sheet CUL, row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 696.6087
$ws.Range("I113").Value = 599.4
$ws.Range("J113").Value = 771.38464
$ws.Range("K113").Value = 1798.2
$ws.Range("L113").Value = 2314.15392
$ws.Range("M113").Value = 371.8000000000002
$ws.Range("N113").Value = -6654.15392

# Sheet CUL, row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1798.8077
$ws.Range("J129").Value = 2769.1538
$ws.Range("L129").Value = 8307.4614
$ws.Range("N129").Value = -18307.4614

# Sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1525.7894
$ws.Range("I132").Value = 1660.7693
$ws.Range("J132").Value = 1233.3334
$ws.Range("K132").Value = 14946.9237
$ws.Range("L132").Value = 11100.0006
$ws.Range("M132").Value = -12416.9237
$ws.Range("N132").Value = -16160.0006

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2479.525
$ws.Range("I122").Value = 2253.9
$ws.Range("K122").Value = 6761.700000000001
$ws.Range("M122").Value = -4311.700000000001

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2282.8206
$ws.Range("I132").Value = 2128.2424
$ws.Range("J132").Value = 3133
$ws.Range("K132").Value = 6384.7272
$ws.Range("L132").Value = 9399
$ws.Range("M132").Value = -3854.7272
$ws.Range("N132").Value = -14459

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 59800.555
$ws.Range("I40").Value = 101480
$ws.Range("J40").Value = 7701.25
$ws.Range("K40").Value = 101480
$ws.Range("L40").Value = 7701.25
$ws.Range("M40").Value = -101344
$ws.Range("N40").Value = -7973.25

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1189
$ws.Range("I61").Value = 1115.6428
$ws.Range("J61").Value = 1702.5
$ws.Range("K61").Value = 1115.6428
$ws.Range("L61").Value = 1702.5
$ws.Range("M61").Value = -913.6428000000001
$ws.Range("N61").Value = -2106.5

# Sheet LTW, row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 11230.7
$ws.Range("I93").Value = 15271.857
$ws.Range("J93").Value = 1801.3334
$ws.Range("K93").Value = 15271.857
$ws.Range("L93").Value = 1801.3334
$ws.Range("M93").Value = -14023.857
$ws.Range("N93").Value = -4297.3334

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1189
$ws.Range("I113").Value = 1115.6428
$ws.Range("J113").Value = 1702.5
$ws.Range("K113").Value = 1115.6428
$ws.Range("L113").Value = 1702.5
$ws.Range("M113").Value = 1054.3572
$ws.Range("N113").Value = -6042.5

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6201.077
$ws.Range("I122").Value = 5060.3
$ws.Range("K122").Value = 15180.9
$ws.Range("M122").Value = -12730.9

# Sheet WVR, row 45
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 18684.75
$ws.Range("J45").Value = 18684.75
$ws.Range("L45").Value = 18684.75
$ws.Range("N45").Value = -19666.75

# Sheet WVR, row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1229.4117
$ws.Range("I81").Value = 840
$ws.Range("J81").Value = 1785.7142
$ws.Range("K81").Value = 1680
$ws.Range("L81").Value = 3571.4284
$ws.Range("M81").Value = -619
$ws.Range("N81").Value = -5693.4284

# Sheet WVR, row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1229.4117
$ws.Range("I84").Value = 840
$ws.Range("J84").Value = 1785.7142
$ws.Range("K84").Value = 8400
$ws.Range("L84").Value = 17857.142
$ws.Range("M84").Value = -3096
$ws.Range("N84").Value = -28465.142

# Sheet WVR, row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1500
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1500
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1500
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -4246

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 435.34482
$ws.Range("I113").Value = 362.35
$ws.Range("J113").Value = 597.55554
$ws.Range("K113").Value = 1087.05
$ws.Range("L113").Value = 1792.66662
$ws.Range("M113").Value = 1082.95
$ws.Range("N113").Value = -6132.66662

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 30455.486
$ws.Range("I122").Value = 33795.29
$ws.Range("J122").Value = 4572
$ws.Range("K122").Value = 101385.87
$ws.Range("L122").Value = 13716
$ws.Range("M122").Value = -98935.87
$ws.Range("N122").Value = -18616

# Sheet WVR, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 51069.75
$ws.Range("I126").Value = 63505.625
$ws.Range("J126").Value = 1326.25
$ws.Range("K126").Value = 190516.875
$ws.Range("L126").Value = 3978.75
$ws.Range("M126").Value = -188046.875
$ws.Range("N126").Value = -8918.75

# Sheet WVR, row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 49991.6
$ws.Range("J135").Value = 49991.6
$ws.Range("L135").Value = 49991.6
$ws.Range("N135").Value = -60131.6

# Sheet WVR, row 137
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 68328.336
$ws.Range("J137").Value = 68328.336
$ws.Range("L137").Value = 68328.336
$ws.Range("N137").Value = -78528.336

# Sheet WVR, row 141
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 50700
$ws.Range("J141").Value = 50700
$ws.Range("L141").Value = 50700
$ws.Range("N141").Value = -61060
